# refactor currency conversion, now explicit source and target amounts
#
# currency_conversions sheet (sheet6): the single "foreign_amount" column is
# replaced by two explicit columns - "source_amount" (new, bold) and
# "target_amount" (reuses the old foreign_amount value) - and a new
# "target_fees" column is added next to the existing "source_fees" column.
# money_transfers (sheet7) becomes the active sheet no longer, instead
# currency_conversions becomes active/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Make currency_conversions the active/selected sheet (was money_transfers).
$ws.Activate()

# --- Header row (row 1) -----------------------------------------------
# Write the brand-new header strings first, and in the exact order that
# they should be appended to the shared string table:
#   target_amount, source_amount, target_fees
$ws.Range("E1").Value = "target_amount"
$ws.Range("B1").Value = "source_amount"
$ws.Range("F1").Value = "target_fees"

# Re-write the remaining (pre-existing) headers at their new positions.
$ws.Range("A1").Value = "date"
$ws.Range("C1").Value = "source_fees"
$ws.Range("D1").Value = "source_currency"
$ws.Range("G1").Value = "target_currency"
$ws.Range("G1").Font.Bold = $true
$ws.Range("H1").Value = "comment"
$ws.Range("H1").Font.Bold = $true

# "source_amount" header gets its own distinct bold style (a fresh font
# entry, as happens when formatting is (re)applied explicitly rather than
# inherited from an existing cell style).
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 12
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.ThemeColor = 1

# --- Data row (row 2) ---------------------------------------------------
# date (A2) is untouched.
$ws.Range("B2").Value = -1       # source_amount
$ws.Range("C2").Value = 0        # source_fees
$ws.Range("D2").Value = "EUR"    # source_currency
$ws.Range("E2").Value = 144.74   # target_amount (was foreign_amount)
$ws.Range("F2").Value = 0        # target_fees
$ws.Range("G2").Value = "USD"    # target_currency
